# Remove bottom logo on slide 7
#
# The "PIT MutationMate" logo picture sitting in the bottom-right corner of
# slide 7 (shape "Google Shape;145;p19", a p:pic referencing media/image1.png)
# is no longer wanted on that slide, so it gets deleted.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Find the picture shape by name (robust against any index shuffling) and
# delete it. Walk backwards so removing a shape never perturbs indices we
# still need to visit.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Google Shape;145;p19") {
        $shape.Delete()
    }
}
